# Update the mandatory-module grade scores for the student in row 15
# (Kamal Badawi, matriculation 5000014) of Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newGrades = @(99, 95, 100, 94, 73, 81, 77, 73, 70, 80, 73, 67, 85, 98, 100, 86, 88, 66, 75, 65, 80, 100, 82, 67, 74, 81)

$grid = New-Object 'object[,]' 1,$newGrades.Length
for ($i = 0; $i -lt $newGrades.Length; $i++) {
    $grid[0,$i] = $newGrades[$i]
}

$ws.Range("E15:AD15").Value = $grid

# Restore the view/selection state recorded for the sheet after the edit.
$ws.Range("AA15").Select()
